# Commit: "Added paths to images in documents"
#
# The author added the absolute path of the "School" image to column P
# ("image") of the data row on Sheet1, and left that new cell selected
# (as Excel does right after you type a value into a cell and press
# Enter/Tab while it remains the active cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New image path value for the "image" column (P) of the data row (row 2).
$ws.Range("P2").Value = "C:Users/vano/Documents/GitHub/ZPI_VAF/iaff_assistant/images/For Parents/school.jpg"

# Leave the freshly-edited cell selected, matching the saved selection state.
[void]$ws.Range("P2").Select()
